# Extend the workbook with a new worksheet "TwoPqAndOnePvNodeDifferentOrde2",
# a copy/extension of "TwoPqAndOnePvNodeDifferentOrder" that also reports the
# real/imaginary parts of the complex numbers and the four "change matrix"
# quadrants (real power by real/imaginary, imaginary power by real/imaginary).

$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("TwoPqAndOnePvNodeDifferentOrder")

# 1. Duplicate the source sheet and move the copy to the end of the workbook.
$src.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "TwoPqAndOnePvNodeDifferentOrde2"

# 2. Add the "real"/"imaginary" headers and IMREAL/IMAGINARY columns next to
#    the existing admittance/magnitude/angle columns (rows 1-4).
$ws.Range("M1").Value = "real"
$ws.Range("Q1").Value = "imaginary"

$ws.Range("M2").Formula = "=IMREAL(A2)"
$ws.Range("N2").Formula = "=IMREAL(B2)"
$ws.Range("O2").Formula = "=IMREAL(C2)"
$ws.Range("M3").Formula = "=IMREAL(A3)"
$ws.Range("N3").Formula = "=IMREAL(B3)"
$ws.Range("O3").Formula = "=IMREAL(C3)"
$ws.Range("M4").Formula = "=IMREAL(A4)"
$ws.Range("N4").Formula = "=IMREAL(B4)"
$ws.Range("O4").Formula = "=IMREAL(C4)"

$ws.Range("Q2").Formula = "=IMAGINARY(A2)"
$ws.Range("R2").Formula = "=IMAGINARY(B2)"
$ws.Range("S2").Formula = "=IMAGINARY(C2)"
$ws.Range("Q3").Formula = "=IMAGINARY(A3)"
$ws.Range("R3").Formula = "=IMAGINARY(B3)"
$ws.Range("S3").Formula = "=IMAGINARY(C3)"
$ws.Range("Q4").Formula = "=IMAGINARY(A4)"
$ws.Range("R4").Formula = "=IMAGINARY(B4)"
$ws.Range("S4").Formula = "=IMAGINARY(C4)"

# 3. Same real/imaginary headers + formulas for the "voltages" block (rows 6-9).
$ws.Range("G6").Value = "real"
$ws.Range("I6").Value = "imaginary"
$ws.Range("G7").Formula = "=IMREAL(A7)"
$ws.Range("I7").Formula = "=IMAGINARY(A7)"
$ws.Range("G8").Formula = "=IMREAL(A8)"
$ws.Range("I8").Formula = "=IMAGINARY(A8)"
$ws.Range("G9").Formula = "=IMREAL(A9)"
$ws.Range("I9").Formula = "=IMAGINARY(A9)"

# 4. Same real/imaginary headers + formulas for the "currents" block (rows 11-14).
$ws.Range("G11").Value = "real"
$ws.Range("I11").Value = "imaginary"
$ws.Range("G12").Formula = "=IMREAL(A12)"
$ws.Range("I12").Formula = "=IMAGINARY(A12)"
$ws.Range("G13").Formula = "=IMREAL(A13)"
$ws.Range("I13").Formula = "=IMAGINARY(A13)"
$ws.Range("G14").Formula = "=IMREAL(A14)"
$ws.Range("I14").Formula = "=IMAGINARY(A14)"

# 5. Replace the old "real power by angle" / "imaginary power by amplitude"
#    blocks (rows 16-23) with the new four change-matrix quadrants.
$ws.Range("A16:M23").ClearContents()

$ws.Range("A16").Value = "real power by real"
$ws.Range("D16").Value = "real power by imaginary"
$ws.Range("G16").Value = "real power by angle"

$ws.Range("A17").Formula = "=(`$N`$2*`$G`$8-`$R`$2*`$I`$8)+(`$O`$2*`$G`$9-`$S`$2*`$I`$9)+2*`$M`$2*`$G`$7-`$G`$12"
$ws.Range("B17").Formula = "=`$G7*N2+`$I7*R2"
$ws.Range("D17").Formula = "=(R2*G8+N2*I8)+(S2*G9+O2*I9)+2*M2*I7-I12"
$ws.Range("E17").Formula = "=`$I7*N2-`$G7*R2"
$ws.Range("G17").Formula = "=C7*G2*C9*SIN(E7-K2-E9)"
$ws.Range("J17").Formula = "=(`$N`$2*`$G`$8-`$R`$2*`$I`$8)"
$ws.Range("K17").Formula = "=`$O`$2*`$G`$9-`$S`$2*`$I`$9"
$ws.Range("L17").Formula = "=2*`$M`$2*`$G`$7"
$ws.Range("M17").Formula = "=-`$G`$12"

$ws.Range("A18").Formula = "=`$G8*M3+`$I8*Q3"
$ws.Range("B18").Formula = "=(M3*G7-Q3*I7)+(O3*G9-S3*I9)+2*N3*G8-G13"
$ws.Range("D18").Formula = "=`$I8*M3-`$G8*Q3"
$ws.Range("E18").Formula = "=(Q3*G7+M3*I7)+(S3*G9+O3*I9)+2*N3*I8-I13"
$ws.Range("G18").Formula = "=C8*G3*C9*SIN(E8-K3-E9)"
$ws.Range("L18").Formula = "=L17+M17+J17"

$ws.Range("A19").Formula = "=`$G9*M4+`$I9*Q4"
$ws.Range("B19").Formula = "=`$G9*N4+`$I9*R4"
$ws.Range("D19").Formula = "=`$I9*M4-`$G9*Q4"
$ws.Range("E19").Formula = "=`$I9*N4-`$G9*R4"
$ws.Range("G19").Formula = "=-(C9*E4*C7*SIN(E9-I4-E7)+C9*F4*C8*SIN(E9-J4-E8))+C14*C9*SIN(E9-E14)"

$ws.Range("A21").Value = "imaginary power by real"
$ws.Range("D21").Value = "imaginary power by imaginary"
$ws.Range("G21").Value = "imaginary power by angle"

$ws.Range("A22").Formula = "=(-R2*G8-N2*I8)+(-S2*G9-O2*I9)-2*Q2*G7+I12"
$ws.Range("B22").Formula = "=`$I7*N2-`$G7*R2"
$ws.Range("D22").Formula = "=(N2*G8-R2*I8)+(O2*G9-S2*I9)-2*Q2*I7-G12"
$ws.Range("E22").Formula = "=-`$G7*N2-`$I7*R2"
$ws.Range("G22").Formula = "=-C7*G2*C9*COS(E7-K2-E9)"

$ws.Range("A23").Formula = "=`$I8*M3-`$G8*Q3"
$ws.Range("B23").Formula = "=(-Q3*G7-M3*I7)+(-S3*G9-O3*I9)-2*R3*G8+I13"
$ws.Range("D23").Formula = "=-`$G8*M3-`$I8*Q3"
$ws.Range("E23").Formula = "=(M3*G7-Q3*I7)+(O3*G9-S3*I9)-2*R3*I8-G13"
$ws.Range("G23").Formula = "=-C8*G3*C9*COS(E8-K3-E9)"

# 6. View state: selections on each sheet, and which sheet/cell is active.
$src.Range("A6:E14").Select()

$ws3 = $wb.Worksheets.Item("ThreePqNodes")
$ws3.Range("H24").Select()

$ws.Activate()
$ws.Range("K20").Select()
